$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.622.74"
$ws.Range("E2").Value = "  -3.27%  "
$ws.Range("D3").Value = "1.850.86"
$ws.Range("E3").Value = "  -3.82%  "
$ws.Range("E4").Value = "  -1.11%  "
$ws.Range("D5").Formula = '="335.69"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("D6").Formula = '="1.002"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  -0.92%  "
$ws.Range("D7").Formula = '="0.4671"'
$ws.Range("D7").Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value = "  -2.92%  "
$ws.Range("D8").Formula = '="0.3904"'
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = "  -3.67%  "
$ws.Range("D9").Formula = '="46.23"'
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = "  -2.66%  "
$ws.Range("D10").Formula = '="0.07911"'
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = "  -3.37%  "
$ws.Range("D12").Formula = '="22.27"'
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = "  -6.29%  "
$ws.Range("D13").Value = "1.907.51"
$ws.Range("E13").Value = "  +3.76%  "
$ws.Range("D14").Formula = '="5.818"'
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = "  -4.16%  "
$ws.Range("D15").Formula = '="6.964"'
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = "  -4.33%  "
$ws.Range("D16").Formula = '="0.06907"'
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value = "  +0.60%  "
$ws.Range("D17").Formula = '="87.76"'
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Value = "  -4.02%  "
$ws.Range("E18").Value = "  -1.22%  "
$ws.Range("E19").Value = "  -3.38%  "
$ws.Range("D20").Formula = '="17.05"'
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = "  -3.15%  "
$ws.Range("D21").Formula = '="1.002"'
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  -0.79%  "
$ws.Range("D22").Value = "28.648.06"
$ws.Range("E22").Value = "  -3.14%  "
$ws.Range("D23").Formula = '="5.391"'
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = "  -4.62%  "
$ws.Range("D24").Formula = '="11.23"'
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  -5.84%  "
$ws.Range("D25").Formula = '="2.157"'
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = "  -1.04%  "
$ws.Range("D26").Value = "2.103.65"
$ws.Range("E26").Value = "  +1.82%  "
$ws.Range("D27").Formula = '="153.29"'
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = "  -1.86%  "
$ws.Range("D28").Formula = '="19.38"'
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = "  -3.22%  "
$ws.Range("D29").Formula = '="6.067"'
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = "  -4.36%  "
$ws.Range("D30").Formula = '="2.014"'
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = "  -3.35%  "
$ws.Range("D31").Formula = '="117.45"'
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = "  -2.59%  "
$ws.Range("D32").Formula = '="0.9680"'
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = "  -3.46%  "
$ws.Range("D33").Formula = '="0.09349"'
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = "  -2.56%  "
$ws.Range("D34").Formula = '="5.365"'
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = "  -4.31%  "
$ws.Range("D35").Formula = '="3.475"'
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = "  -2.32%  "
$ws.Range("D36").Formula = '="1.346"'
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = "  -3.12%  "
$ws.Range("D37").Formula = '="0.06116"'
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = "  -5.85%  "
$ws.Range("D38").Formula = '="0.02205"'
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = "  -3.07%  "
$ws.Range("D39").Formula = '="1.170"'
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = "  -3.80%  "
$ws.Range("D40").Formula = '="7.681"'
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = "  -2.10%  "
$ws.Range("D41").Formula = '="0.5698"'
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = "  -3.92%  "
$ws.Range("D42").Formula = '="10.11"'
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = "  -5.75%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Formula = '="0.1791"'
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = "  -2.81%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Formula = '="2.428"'
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = "  -2.63%  "
$ws.Range("D45").Formula = '="1.249"'
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = "  -2.47%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Formula = '="11.80"'
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = "  -4.76%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Formula = '="0.5372"'
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = "  -3.06%  "
$ws.Range("D48").Formula = '="0.07099"'
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = "  -5.76%  "
$ws.Range("D49").Formula = '="1.906"'
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = "  -2.26%  "
$ws.Range("D50").Formula = '="113.06"'
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = "  -4.48%  "
$ws.Range("D51").Formula = '="2.353"'
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = "  -3.11%  "
$ws.Range("A1").Select()
